# "Added clean you can be God's Friend"
# Appends a new "Were Mito Ni Ibedi Mere Pere" (You Can Be God's Friend) Q&A
# section to the bottom of the "Home Page" sheet, and nudges the saved
# selection on the "Common Verbs" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Body rows (wrap-text, single column) + the new section title (row 99).
# Values are written top-to-bottom so the shared-string table picks up
# the same append order as the authored workbook. Row 96 is left blank
# on purpose as a spacer before the new section.
# ---------------------------------------------------------------------
$ws1.Cells.Item(97, 1).Value = 'Were Mito Ni Ibedi Mere Pere'
$ws1.Cells.Item(97, 1).WrapText = $true

$ws1.Cells.Item(98, 1).Value = 'Go lakonyin limo kisangala i kwo.'
$ws1.Cells.Item(98, 1).WrapText = $true

# Row 99 title: reuse the formatting already used for the other section
# heading in this workbook (sheet2 row 29: Arial, 12pt, #222222) and only
# recolor it, so exactly one new font / cellXf pair is produced instead
# of a pile of intermediate ones.
$ws2.Cells.Item(29, 1).Copy() | Out-Null
$ws1.Cells.Item(99, 1).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(99, 1).Font.Color = 2236962
$ws1.Rows.Item(99).RowHeight = 15.75
$ws1.Cells.Item(99, 1).Value = 'Ginaŋo manyalo konyin ŋeyo gima go mito kodi gima go kimiti.'

$ws1.Cells.Item(100, 1).Value = 'Go konyo wan ŋeyo gima go otimo cango con, gima go tima pama, kodi gima go latimo i hongo mabino.'
$ws1.Cells.Item(100, 1).WrapText = $true
$ws1.Rows.Item(100).RowHeight = 30

$ws1.Cells.Item(101, 1).Value = 'Piny manyien k’obedi paka piny ma wan’iye pama. Piny manyien ochale nedi?'
$ws1.Cells.Item(101, 1).WrapText = $true
$ws1.Rows.Item(101).RowHeight = 30

$ws1.Cells.Item(102, 1).Value = 'Aŋo manyutho?'
$ws1.Cells.Item(102, 1).WrapText = $true

$ws1.Cells.Item(103, 1).Value = 'Gima Baibul waco ri wan kwoŋ mikula ma ndir pa Nowa fonjo wan aŋo?'
$ws1.Cells.Item(103, 1).WrapText = $true
$ws1.Rows.Item(103).RowHeight = 30

# ---------------------------------------------------------------------
# Printable page setup for the Home Page sheet (matches Common Verbs).
# ---------------------------------------------------------------------
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Saved selections: Home Page lands on the new section (B97), Common
# Verbs' cursor moves to A45 - but keep Home Page as the active tab.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A45").Select() | Out-Null

$ws1.Activate()
$ws1.Range("B97").Select() | Out-Null
